$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '310.29'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = '-1.03%'
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '37.56'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = '-4.19%'
$ws.Range("E3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '5.085'
$ws.Range("D4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '0.07762'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = '-4.80%'
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '4.345'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = '-3.07%'
$ws.Range("E6").ClearFormats()

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value2 = 'KuCoinToken'
$ws.Range("B7").ClearFormats()

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value2 = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("C7").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '8.207'
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = '-1.20%'
$ws.Range("E7").ClearFormats()

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value2 = 'FTXToken'
$ws.Range("B8").ClearFormats()

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value2 = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("C8").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '1.887'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = '-3.92%'
$ws.Range("E8").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = '-10.88%'
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.9156'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = '-2.68%'
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.1201'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = '-8.88%'
$ws.Range("E11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '0.1919'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = '-2.44%'
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '0.08908'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = '-1.10%'
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '0.03420'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = '-2.23%'
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '0.09694'
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = '-0.24%'
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '0.001369'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = '-2.78%'
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '0.005911'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = '-3.80%'
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '3.552'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = '-0.71%'
$ws.Range("E18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '0.3409'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = '-1.63%'
$ws.Range("E19").ClearFormats()

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value2 = 'ProBitToken'
$ws.Range("B20").ClearFormats()

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value2 = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("C20").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '0.1279'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = '-1.74%'
$ws.Range("E20").ClearFormats()

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value2 = 'MCDex'
$ws.Range("B21").ClearFormats()

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value2 = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("C21").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '5.033'
$ws.Range("D21").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = '0.36%'
$ws.Range("E21").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = '3.82%'
$ws.Range("E22").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = '5,586.25%'
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '0.04377'
$ws.Range("D24").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = '-0.06%'
$ws.Range("E24").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = '-2.86%'
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '0.004248'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = '-10.18%'
$ws.Range("E26").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '0.02108'
$ws.Range("D39").ClearFormats()

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = '-4.58%'
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '0.04945'
$ws.Range("D40").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = '-5.40%'
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '0.007658'
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = '1.23%'
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '0.009874'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = '-4.32%'
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.1341'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = '-3.77%'
$ws.Range("E43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.002059'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = '-2.18%'
$ws.Range("E44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '0.009584'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = '5.08%'
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '0.00006705'
$ws.Range("D46").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = '-1.79%'
$ws.Range("E46").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = '-0.27%'
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '0.003041'
$ws.Range("D48").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = '0.74%'
$ws.Range("E48").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = '-0.27%'
$ws.Range("E50").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = '-0.27%'
$ws.Range("E51").ClearFormats()
